# Ontario daily doses — append two more days of data (rows 98 & 99),
# extending the running-average / cumulative formulas down to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 98 (record #96, 2021-04-04) -----------------------------------
$ws.Range("A98").Value = 96
$ws.Range("B98").Value = 44290
$ws.Range("C98").Value = 69125
$ws.Range("D98").Value = 2493188
$ws.Range("E98").Value = 644095
$ws.Range("F98").Value = 322048

# --- New row 99 (record #97, 2021-04-05) -----------------------------------
$ws.Range("A99").Value = 97
$ws.Range("B99").Value = 44291
$ws.Range("C99").Value = 52452
$ws.Range("D99").Value = 2545640
$ws.Range("E99").Value = 644393
$ws.Range("F99").Value = 322197

# Match the existing "#,##0"-style numeric formatting used throughout C:I.
$ws.Range("C98:F99").NumberFormat = "#,##0"

# Extend the 7-day-average column. Filling G91:G99 in one go merges the
# previously-individual G91:G97 formulas into a single shared group and
# continues it down through the two new rows.
$ws.Range("G91:G97").Formula = "=AVERAGE(C85:C91)"
$ws.Range("G98:G99").Formula = "=AVERAGE(C92:C98)"

# Extend the running "doses left" and "doses administered vs population"
# columns down through the new rows as their own shared-formula group.
$ws.Range("H98:H99").Formula = "=H97-C98"
$ws.Range("I98:I99").Formula = "=(D98-F98)"

# Reflect the scrolled/selected state shown in the saved workbook.
$ws.Range("H103").Select()
